# Add an "Autre identifiant" line to the etablissement block in the
# TIAC/SSA "evenement simple" template, right after the "N SIRET" line
# and before the "Raison sociale" line:
#
#   Autre identifiant : {{ etablissement.autre_identifiant }}

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# Locate the paragraph that renders the SIRET field without touching/
# re-walking every paragraph in the document.
$rng = $d.Content
[void]$rng.Find.Execute("etablissement.siret", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$siretPara = $rng.Paragraphs(1)

# Create a new (empty) paragraph right after it; Word copies the
# paragraph/run formatting (style "Corpsdetexte", Calibri 10pt) from the
# paragraph it was split from, which matches the source paragraphs
# around it.
$siretPara.Range.InsertParagraphAfter()

# Build the sentence as four separate runs (mirrors how the equivalent
# "Numero d'inspection Resytal" field further down is split into
# multiple runs around the templated field name) by typing each chunk
# into its own paragraph, then joining the paragraphs back together by
# deleting the paragraph marks between them. This keeps each chunk as
# its own <w:r> instead of Word silently coalescing same-format runs
# when text is simply appended in place.
$p1Start = $siretPara.Range.End
$p1 = $d.Range($p1Start, $p1Start).Paragraphs(1)
$p1.Range.InsertAfter("Autre identifiant")

$p1.Range.InsertParagraphAfter()
$p2Start = $p1.Range.End
$p2 = $d.Range($p2Start, $p2Start).Paragraphs(1)
$p2.Range.InsertAfter("${nbsp}: {{ etablissement.")

$p2.Range.InsertParagraphAfter()
$p3Start = $p2.Range.End
$p3 = $d.Range($p3Start, $p3Start).Paragraphs(1)
$p3.Range.InsertAfter("autre_identifiant")

$p3.Range.InsertParagraphAfter()
$p4Start = $p3.Range.End
$p4 = $d.Range($p4Start, $p4Start).Paragraphs(1)
$p4.Range.InsertAfter(" }}")

# Join p1..p4 back into a single paragraph (delete the 3 paragraph
# marks that currently separate them).
$e = $p1.Range.End
$d.Range($e - 1, $e).Delete()
$e = $p1.Range.End
$d.Range($e - 1, $e).Delete()
$e = $p1.Range.End
$d.Range($e - 1, $e).Delete()

Write-Output $p1.Range.Text
